$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 609, shifting rows 609:651 down to 610:652
$ws.Rows("609:609").Insert(-4121)  # xlShiftDown

# Populate the new row 609 with the new data record
$ws.Range("A609").Value = 3
$ws.Range("B609").Value = "Femacal de La Calera"
$ws.Range("C609").Value = "Coquimbo"
$ws.Range("D609").Value = 45265
$ws.Range("D609").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E609").Value = 5
$ws.Range("F609").Value = 100112012
$ws.Range("G609").Value = "Espinaca"
$ws.Range("H609").Value = "Sin especificar"
$ws.Range("I609").Value = "Primera"
$ws.Range("J609").Value = 80
$ws.Range("K609").Value = 5000
$ws.Range("L609").Value = 5000
$ws.Range("M609").Value = 5000
$ws.Range("N609").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O609").Value = "Provincia de Quillota"
$ws.Range("P609").Value = 1667
$ws.Range("Q609").Value = 3
$ws.Range("R609").Value = "Hortaliza"
